# Edit: use hl60 RNA-seq 72h data instead of 120h.
# This drops the "Stbl.*" (120h stability) columns entirely and updates the
# "Exp.hl60.log2FC" column (and the downstream TE.* columns) with the 72h values,
# which also changes the row order of the "Comparisons" and "CRISPRi-screens" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Comparisons": drop the Stbl.* columns (P:AA) then rewrite rows
# 2-9 in the new gene order with the updated 72h values.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Comparisons")
$ws2.Range("P1:AA9").EntireColumn.Delete()

$ws2.Range("A2").Value = "ELP5"
$ws2.Range("B2").Value = 1.072
$ws2.Range("C2").Value = 0.005
$ws2.Range("D2").Value = 0.075
$ws2.Range("E2").Value = 0.766
$ws2.Range("F2").Value = -0.221
$ws2.Range("G2").Value = 0.37
$ws2.Range("H2").Value = 0.259
$ws2.Range("I2").Value = 0.397
$ws2.Range("J2").Value = 0.2
$ws2.Range("K2").Value = 0.394
$ws2.Range("L2").Value = 0.782
$ws2.Range("M2").Value = 0.014
$ws2.Range("N2").Value = -0.047
$ws2.Range("O2").Value = 0.978

$ws2.Range("A3").Value = "KLF16"
$ws2.Range("B3").Value = 0.962
$ws2.Range("C3").Value = 0.001
$ws2.Range("D3").Value = 0.048
$ws2.Range("E3").Value = 0.858
$ws2.Range("F3").Value = -0.354
$ws2.Range("G3").Value = 0.238
$ws2.Range("H3").Value = 0.113
$ws2.Range("I3").Value = 0.557
$ws2.Range("J3").Value = -0.36
$ws2.Range("K3").Value = 0.308
$ws2.Range("L3").Value = -0.433
$ws2.Range("M3").Value = 0.111
$ws2.Range("N3").Value = -0.183
$ws2.Range("O3").Value = 0.91

$ws2.Range("A4").Value = "FDFT1"
$ws2.Range("B4").Value = 0.788
$ws2.Range("C4").Value = 0.001
$ws2.Range("D4").Value = 0.05
$ws2.Range("E4").Value = 0.729
$ws2.Range("F4").Value = 0.202
$ws2.Range("G4").Value = 0.35
$ws2.Range("H4").Value = 0.099
$ws2.Range("I4").Value = 0.633
$ws2.Range("J4").Value = -0.288
$ws2.Range("K4").Value = 0.221
$ws2.Range("L4").Value = -0.076
$ws2.Range("M4").Value = 0.73
$ws2.Range("N4").Value = 0.215
$ws2.Range("O4").Value = 0.464

$ws2.Range("A5").Value = "MECR"
$ws2.Range("B5").Value = 1.6
$ws2.Range("C5").Value = 0.006
$ws2.Range("D5").Value = -0.184
$ws2.Range("E5").Value = 0.423
$ws2.Range("F5").Value = -0.523
$ws2.Range("G5").Value = 0.055
$ws2.Range("H5").Value = -0.198
$ws2.Range("I5").Value = 0.441
$ws2.Range("J5").Value = -0.776
$ws2.Range("K5").Value = 0.002
$ws2.Range("L5").Value = 0.198
$ws2.Range("M5").Value = 0.632
$ws2.Range("N5").Value = -0.044
$ws2.Range("O5").Value = 0.985

$ws2.Range("A6").Value = "MIOS"
$ws2.Range("B6").Value = 0.626
$ws2.Range("C6").Value = 0.002
$ws2.Range("D6").Value = -0.133
$ws2.Range("E6").Value = 0.492
$ws2.Range("F6").Value = -0.141
$ws2.Range("G6").Value = 0.513
$ws2.Range("H6").Value = -0.047
$ws2.Range("I6").Value = 0.807
$ws2.Range("J6").Value = 0.293
$ws2.Range("K6").Value = 0.146
$ws2.Range("L6").Value = 0.373
$ws2.Range("M6").Value = 0.066
$ws2.Range("N6").Value = -0.26
$ws2.Range("O6").Value = 0.842

$ws2.Range("A7").Value = "TTI1"
$ws2.Range("B7").Value = 0.537
$ws2.Range("C7").Value = 0.011
$ws2.Range("D7").Value = -0.214
$ws2.Range("E7").Value = 0.376
$ws2.Range("F7").Value = -0.422
$ws2.Range("G7").Value = 0.073
$ws2.Range("H7").Value = -0.331
$ws2.Range("I7").Value = 0.095
$ws2.Range("J7").Value = -0.184
$ws2.Range("K7").Value = 0.338
$ws2.Range("L7").Value = -0.016
$ws2.Range("M7").Value = 0.933
$ws2.Range("N7").ClearContents()
$ws2.Range("O7").ClearContents()

$ws2.Range("A8").Value = "NUDCD3"
$ws2.Range("B8").Value = 0.33
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = -0.313
$ws2.Range("E8").Value = 0.051
$ws2.Range("F8").Value = -0.466
$ws2.Range("G8").Value = 0.01
$ws2.Range("H8").Value = -0.054
$ws2.Range("I8").Value = 0.769
$ws2.Range("J8").Value = 0.004
$ws2.Range("K8").Value = 0.986
$ws2.Range("L8").Value = -0.238
$ws2.Range("M8").Value = 0.27
$ws2.Range("N8").Value = 0.698
$ws2.Range("O8").Value = 0.333

$ws2.Range("A9").Value = "ZNF787"
$ws2.Range("B9").Value = 1.306
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = -0.394
$ws2.Range("E9").Value = 0.108
$ws2.Range("F9").Value = -0.138
$ws2.Range("G9").Value = 0.664
$ws2.Range("H9").Value = -0.084
$ws2.Range("I9").Value = 0.701
$ws2.Range("J9").Value = -0.279
$ws2.Range("K9").Value = 0.267
$ws2.Range("L9").Value = -0.201
$ws2.Range("M9").Value = 0.347
$ws2.Range("N9").Value = 0.504
$ws2.Range("O9").Value = 0.745

# ---------------------------------------------------------------------
# Sheet "CRISPRi-screens": same gene reordering (values unchanged).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CRISPRi-screens")

$ws3.Range("A2").Value = "ELP5"
$ws3.Range("B2").Value = 0.233
$ws3.Range("C2").Value = 0.004
$ws3.Range("D2").Value = 0.136
$ws3.Range("E2").Value = 0.083
$ws3.Range("F2").Value = 0.078
$ws3.Range("G2").Value = 0.311
$ws3.Range("H2").Value = -0.191
$ws3.Range("I2").Value = 0.442
$ws3.Range("J2").Value = -0.489
$ws3.Range("K2").Value = 0.06
$ws3.Range("L2").Value = -0.489
$ws3.Range("M2").Value = 0.06

$ws3.Range("A3").Value = "KLF16"
$ws3.Range("B3").Value = 0.165
$ws3.Range("C3").Value = 0.005
$ws3.Range("D3").Value = 0.204
$ws3.Range("E3").Value = 0.005
$ws3.Range("F3").Value = 0.283
$ws3.Range("G3").Value = 0.02
$ws3.Range("H3").Value = -0.069
$ws3.Range("I3").Value = 0.001
$ws3.Range("J3").Value = -0.108
$ws3.Range("K3").Value = 0.103
$ws3.Range("L3").Value = -0.108
$ws3.Range("M3").Value = 0.103

$ws3.Range("A4").Value = "FDFT1"
$ws3.Range("B4").Value = 0.124
$ws3.Range("C4").Value = 0.006
$ws3.Range("D4").Value = 0.2
$ws3.Range("E4").Value = 0.01
$ws3.Range("F4").Value = 0.203
$ws3.Range("G4").Value = 0.016
$ws3.Range("H4").Value = -0.151
$ws3.Range("I4").Value = 0
$ws3.Range("J4").Value = -0.218
$ws3.Range("K4").Value = 0.001
$ws3.Range("L4").Value = -0.218
$ws3.Range("M4").Value = 0.001

$ws3.Range("A5").Value = "MECR"
$ws3.Range("B5").Value = 0.104
$ws3.Range("C5").Value = 0.234
$ws3.Range("D5").Value = 0.126
$ws3.Range("E5").Value = 0.244
$ws3.Range("F5").Value = 0.245
$ws3.Range("G5").Value = 0.124
$ws3.Range("H5").Value = -0.14
$ws3.Range("I5").Value = 0.096
$ws3.Range("J5").Value = -0.178
$ws3.Range("K5").Value = 0.115
$ws3.Range("L5").Value = -0.178
$ws3.Range("M5").Value = 0.115

$ws3.Range("A6").Value = "MIOS"
$ws3.Range("B6").Value = 0.148
$ws3.Range("C6").Value = 0
$ws3.Range("D6").Value = 0.317
$ws3.Range("E6").Value = 0.004
$ws3.Range("F6").Value = 0.414
$ws3.Range("G6").Value = 0.001
$ws3.Range("H6").Value = -0.114
$ws3.Range("I6").Value = 0.001
$ws3.Range("J6").Value = -0.196
$ws3.Range("K6").Value = 0.001
$ws3.Range("L6").Value = -0.196
$ws3.Range("M6").Value = 0.001

$ws3.Range("A7").Value = "TTI1"
$ws3.Range("B7").Value = 0.192
$ws3.Range("C7").Value = 0.012
$ws3.Range("D7").Value = 0.244
$ws3.Range("E7").Value = 0.213
$ws3.Range("F7").Value = 0.259
$ws3.Range("G7").Value = 0.099
$ws3.Range("H7").Value = -0.4
$ws3.Range("I7").Value = 0.003
$ws3.Range("J7").Value = -0.296
$ws3.Range("K7").Value = 0.044
$ws3.Range("L7").Value = -0.296
$ws3.Range("M7").Value = 0.044

$ws3.Range("A8").Value = "NUDCD3"
$ws3.Range("B8").Value = 0.11
$ws3.Range("C8").Value = 0.134
$ws3.Range("D8").Value = 0.186
$ws3.Range("E8").Value = 0.135
$ws3.Range("F8").Value = 0.382
$ws3.Range("G8").Value = 0.087
$ws3.Range("H8").Value = -0.476
$ws3.Range("I8").Value = 0
$ws3.Range("J8").Value = -0.829
$ws3.Range("K8").Value = 0
$ws3.Range("L8").Value = -0.829
$ws3.Range("M8").Value = 0

$ws3.Range("A9").Value = "ZNF787"
$ws3.Range("B9").Value = 0.114
$ws3.Range("C9").Value = 0.046
$ws3.Range("D9").Value = 0.244
$ws3.Range("E9").Value = 0.039
$ws3.Range("F9").Value = 0.246
$ws3.Range("G9").Value = 0.117
$ws3.Range("H9").Value = -0.068
$ws3.Range("I9").Value = 0.001
$ws3.Range("J9").Value = -0.089
$ws3.Range("K9").Value = 0.485
$ws3.Range("L9").Value = -0.089
$ws3.Range("M9").Value = 0.485
